$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 6 data
$ws.Range("A6").Value = 43669
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

$ws.Range("B6").Formula = "=C5"
$ws.Range("C6").Value = 5851
$ws.Range("D6").Formula = "=C6-B6"
$ws.Range("E6").Value = 0.75

# Move the active selection to F6
$ws.Range("F6").Select()
